$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New paragraph marker rows inserted at row 5 and row 7 (they repeat the
# "paragraph" / "style=smack my ass" header pair used at row 2), the old
# "Lorem ipsum..." content moves down to row 8, and row 6 becomes a fully
# empty row (no values, no styles) representing a preserved empty paragraph.

$ws.Range("A5").Font.Name = "Source Code Pro"
$ws.Range("A5").Value = "paragraph"
$ws.Range("B5").Value = "style=smack my ass"

$ws.Range("A6").Clear()
$ws.Range("B6").Clear()

$ws.Range("A7").Font.Name = "Source Code Pro"
$ws.Range("A7").Value = "paragraph"
$ws.Range("B7").Value = "style=smack my ass"

$ws.Range("B8").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum."

$ws.Range("C6").Select()
